# This script applies a "rotation" of the weekly price records held in rows 2-8.
# The last two rows (7 and 8) of the existing data move up to become the new
# rows 2 and 3, while the old rows 2-6 shift down to become the new rows 4-8.
# Columns A, B, C, E, F, G, H, O, R are identical across these rows and are
# left untouched; only D, I, J, K, L, M, N, P, Q need to be rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the columns that change, for source rows 2-8.
$D = @{}
$I = @{}
$J = @{}
$K = @{}
$L = @{}
$M = @{}
$N = @{}
$P = @{}
$Q = @{}

for ($r = 2; $r -le 8; $r++) {
    $D[$r] = $ws.Cells.Item($r, 4).Value2
    $I[$r] = $ws.Cells.Item($r, 9).Value2
    $J[$r] = $ws.Cells.Item($r, 10).Value2
    $K[$r] = $ws.Cells.Item($r, 11).Value2
    $L[$r] = $ws.Cells.Item($r, 12).Value2
    $M[$r] = $ws.Cells.Item($r, 13).Value2
    $N[$r] = $ws.Cells.Item($r, 14).Value2
    $P[$r] = $ws.Cells.Item($r, 16).Value2
    $Q[$r] = $ws.Cells.Item($r, 17).Value2
}

# Mapping of new row -> old (source) row, implementing the rotation.
$mapping = @{
    2 = 7
    3 = 8
    4 = 2
    5 = 3
    6 = 4
    7 = 5
    8 = 6
}

foreach ($newRow in 2..8) {
    $srcRow = $mapping[$newRow]

    $ws.Cells.Item($newRow, 4).Value2 = $D[$srcRow]
    $ws.Cells.Item($newRow, 9).Value2 = $I[$srcRow]
    $ws.Cells.Item($newRow, 10).Value2 = $J[$srcRow]
    $ws.Cells.Item($newRow, 11).Value2 = $K[$srcRow]
    $ws.Cells.Item($newRow, 12).Value2 = $L[$srcRow]
    $ws.Cells.Item($newRow, 13).Value2 = $M[$srcRow]
    $ws.Cells.Item($newRow, 14).Value2 = $N[$srcRow]
    $ws.Cells.Item($newRow, 16).Value2 = $P[$srcRow]
    $ws.Cells.Item($newRow, 17).Value2 = $Q[$srcRow]
}
